$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.198958039283752
$ws.Range("B1").Value = 1.704150795936584
$ws.Range("C1").Value = 4.635594367980957
$ws.Range("D1").Value = 0.71562659740448
$ws.Range("E1").Value = 0.7683822512626648
